$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update A2 value
$ws.Range("A2").Value = 2000033005

# Fill A3 with value and switch its style to "center horizontal only" (matches
# style index 2 used elsewhere in the sheet, e.g. the header cell A1)
$ws.Range("A3").Value = 2000033006
$ws.Range("A3").Style = "Normal"
$ws.Range("A3").HorizontalAlignment = -4108  # xlCenter

# Add new A4 with value and the same center-horizontal style
$ws.Range("A4").Value = 2000062283
$ws.Range("A4").Style = "Normal"
$ws.Range("A4").HorizontalAlignment = -4108  # xlCenter

# Select A4 as the active cell, matching the final selection state
$ws.Range("A4").Select()
